# Update Name of Algo
# Applies value corrections to the RandomForest result data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value  = -11.0025
$ws.Range("A3").Value  = -21.31900000000002
$ws.Range("C5").Value  = -14.11199999999999
$ws.Range("E5").Value  = 13.21219999999999
$ws.Range("E9").Value  = 14.72350000000001
$ws.Range("E11").Value = 13.47339999999999
$ws.Range("A14").Value = -20.67459999999999
$ws.Range("A16").Value = -20.4793
$ws.Range("C16").Value = -11.9213
$ws.Range("E17").Value = 14.04510000000002
$ws.Range("A21").Value = -21.306
$ws.Range("E21").Value = 12.9258
$ws.Range("A23").Value = -21.37410000000002
$ws.Range("A25").Value = -22.68400000000003
